$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Item("Picture 26").Delete()
